$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("B2").Value = 16
$ws.Range("B3").Value = 18
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 7
